$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing input: Direction(deg) 90 -> -10 ---
$ws.Range("D6").Value = -10

# --- Column width updates ---
# (ColumnWidth values chosen so the serialized "width" attribute lands on
# the closest value the engine's character-width grid can represent.)
$ws.Columns.Item(3).ColumnWidth = 18.166666666666668
$ws.Columns.Item(4).ColumnWidth = 12.833333333333334
$ws.Columns.Item(5).ColumnWidth = 12.833333333333334
$ws.Columns.Item(6).ColumnWidth = 11.833333333333334

# --- New "Joystick calculations" section (rows 20-29) ---
# Write cells in the order needed so new shared strings are appended
# in the same order as the target workbook.
$ws.Range("C21").Value = "x"
$ws.Range("C22").Value = "y"
$ws.Range("C23").Value = "angle(rad)"
$ws.Range("F23").Value = "Robot centric"
$ws.Range("C20").Value = "Joystick calculations"
$ws.Range("C25").Value = "Stick speed"
$ws.Range("C26").Value = "Stick normalization"
$ws.Range("F24").Value = "Normalizers"
$ws.Range("C29").Value = "Final speed"

$ws.Range("D21").Value = 1
$ws.Range("D22").Value = 1

$ws.Range("D23").Formula = "=ATAN2(D21,D22)"
$ws.Range("E23").Formula = "=DEGREES(D23)"
$ws.Range("G23").Formula = "=90-E23"
$ws.Range("G23").Interior.Color = $ws.Range("D6").Interior.Color

$ws.Range("G24").Formula = "=MOD(ABS(G23), 90)"
$ws.Range("H24").Formula = "=IF(G24>45, 90-G24, G24)"
$ws.Range("I24").Formula = "=H24*PI()/180"

$ws.Range("D25").Formula = "=SQRT((D21 * D21) + (D22 * D22))"

$ws.Range("D26").Formula = "=SQRT(1+TAN(I24))/SQRT(2)"

$ws.Range("D29").Formula = "=D25/(D26*SQRT(2))"

$ws.Range("D29").Select()
